$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.582.38'
$ws.Range("E2").Value = '  +5.03%  '
$ws.Range("D3").Value = '2.722.36'
$ws.Range("E3").Value = '  +3.67%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.09'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.23'
$ws.Range("E6").Value = '  +5.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("E8").Value = '  +1.51%  '
$ws.Range("D9").Value = '2.747.15'
$ws.Range("E9").Value = '  +4.20%  '
$ws.Range("E10").Value = '  +3.10%  '
$ws.Range("E11").Value = '  +6.22%  '
$ws.Range("E12").Value = '  +3.91%  '
$ws.Range("D14").Value = '3.206.32'
$ws.Range("E14").Value = '  +3.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.30'
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").Value = '63.425.88'
$ws.Range("E16").Value = '  +4.81%  '
$ws.Range("E17").Value = '  +7.26%  '
$ws.Range("D18").Value = '2.740.62'
$ws.Range("E18").Value = '  +4.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.97'
$ws.Range("E19").Value = '  +3.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.87'
$ws.Range("E20").Value = '  +3.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '360.31'
$ws.Range("E21").Value = '  +3.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.96'
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.538'
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.995'
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.91'
$ws.Range("E25").Value = '  +3.33%  '
$ws.Range("E26").Value = '  +4.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.57'
$ws.Range("E27").Value = '  +5.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '0.0₃0901'
$ws.Range("E29").Value = '  +13.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.01'
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.10'
$ws.Range("E31").Value = '  +6.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '172.57'
$ws.Range("E32").Value = '  +1.91%  '
$ws.Range("E33").Value = '  +13.42%  '
$ws.Range("E34").Value = '  -0.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.47'
$ws.Range("E35").Value = '  +4.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.76'
$ws.Range("E36").Value = '  +7.15%  '
$ws.Range("E37").Value = '  +9.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.81'
$ws.Range("E38").Value = '  +9.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.994'
$ws.Range("E39").Value = '  +15.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '343.53'
$ws.Range("E40").Value = '  +3.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.21'
$ws.Range("E41").Value = '  +5.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.10'
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.52'
$ws.Range("E43").Value = '  +6.91%  '
$ws.Range("E44").Value = '  +8.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.80'
$ws.Range("E45").Value = '  +5.73%  '
$ws.Range("E46").Value = '  +6.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '138.93'
$ws.Range("E47").Value = '  +4.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0587'
$ws.Range("E48").Value = '  +5.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0256'
$ws.Range("E49").Value = '  +4.62%  '
$ws.Range("E50").Value = '  +1.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.995'
$ws.Range("E51").Value = '  -0.37%  '
